# Fix Supreme Court mention missing from short data engineering resumes.
#
# Rewrites the four bullet paragraphs under "KEY ACHIEVEMENTS AND IMPACT" >
# "Impact" into six bullet paragraphs, preserving the bold/colored
# percentage & dollar runs, and adding the Supreme Court achievements.

$d = $word.ActiveDocument

$ACCENT_COLOR = 5258796   # BGR int equivalent of w:color val="2C3E50"
$BULLET = [char]0x2022

# --- helpers -------------------------------------------------------------
# Clear a paragraph's content (keep the paragraph mark) and return a
# collapsed Range positioned at the (now empty) start of the paragraph.
function Clear-ParaContent($para) {
    $full = $para.Range
    $r = $d.Range($full.Start, $full.End - 1)
    if ($r.Start -lt $r.End) {
        $r.Delete()
    }
    return $d.Range($full.Start, $full.Start)
}

# Insert plain (unformatted) text right after the collapsed range $pos;
# return a new collapsed range positioned right after the inserted text.
function Add-PlainRun($pos, [string]$text) {
    $pos.InsertAfter($text)
    return $d.Range($pos.End, $pos.End)
}

# Insert bold/colored text right after the collapsed range $pos; return a
# new collapsed range positioned right after the inserted text.
function Add-BoldRun($pos, [string]$text) {
    $pos.InsertAfter($text)
    $ins = $d.Range($pos.Start, $pos.End)
    $ins.Font.Bold = 1
    $ins.Font.Color = $ACCENT_COLOR
    return $d.Range($pos.End, $pos.End)
}

# --- locate the four target bullet paragraphs ----------------------------
# Anchor off the unique "KEY ACHIEVEMENTS AND IMPACT" heading, which is
# immediately followed by the "Impact" sub-heading and then the four
# bullet paragraphs we need to rewrite. This is more robust than a bare
# hard-coded paragraph index.
$headingIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.Trim()
    if ($t -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $headingIdx = $i
        break
    }
}
if ($headingIdx -eq -1) {
    throw "Could not locate 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# headingIdx+1 == "Impact" sub-heading
# 63: "Discovered systematic race coding errors..."
# 64: "Algorithm reduced mapping costs by 73.5%, saving campaigns..."
# 65: "Built redistricting platform used by thousands of analysts..."
# 66: "Achieved 87% prediction accuracy..."
$P1 = $headingIdx + 2
$P2 = $headingIdx + 3
$P3 = $headingIdx + 4
$P4 = $headingIdx + 5

# Sanity-check the anchors actually hold the text we expect before mutating.
$chk1 = $d.Paragraphs($P1).Range.Text.Trim()
$chk2 = $d.Paragraphs($P2).Range.Text.Trim()
$chk3 = $d.Paragraphs($P3).Range.Text.Trim()
$chk4 = $d.Paragraphs($P4).Range.Text.Trim()
if ($chk1 -ne "• Discovered systematic race coding errors affecting all Black and Asian-American voters") {
    throw "Unexpected paragraph at P1: $chk1"
}
if ($chk3 -notlike "*Built redistricting platform used by thousands of analysts nationwide*") {
    throw "Unexpected paragraph at P3: $chk3"
}
if ($chk4 -notlike "*Achieved*prediction accuracy for voter turnout*") {
    throw "Unexpected paragraph at P4: $chk4"
}

# --- Paragraph 1: Algorithmic innovation / 73.5% --------------------------
$pos = Clear-ParaContent $d.Paragraphs($P1)
$t1 = $BULLET + " Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **"
$pos = Add-PlainRun $pos $t1
$pos = Add-BoldRun $pos "73.5%"
$pos = Add-PlainRun $pos "**"

# --- Paragraph 2: $4.7M savings enabled nonprofit access ------------------
$pos = Clear-ParaContent $d.Paragraphs($P2)
$t2 = $BULLET + " **"
$pos = Add-PlainRun $pos $t2
$dollarAmt = "`$4.7M"
$pos = Add-BoldRun $pos $dollarAmt
$t2b = "** savings enabled nonprofit access"
$pos = Add-PlainRun $pos $t2b

# --- Insert 3 new empty paragraphs right before paragraph $P3 -------------
$target = $d.Paragraphs($P3).Range
$insertPoint = $d.Range($target.Start, $target.Start)
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()

# The 3 blank paragraphs now occupy indices P3, P3+1, P3+2; the old
# "Built redistricting..." bullet has shifted to P3+3.
$legalText = $BULLET + " Legal precedent: Data analysis utilized in Supreme Court case"
$expertText = $BULLET + " Expert methodology validated at highest judicial level"
$discoveryText = $BULLET + " Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

$pos = Clear-ParaContent $d.Paragraphs($P3)
$pos = Add-PlainRun $pos $legalText

$pos = Clear-ParaContent $d.Paragraphs($P3 + 1)
$pos = Add-PlainRun $pos $expertText

$pos = Clear-ParaContent $d.Paragraphs($P3 + 2)
$pos = Add-PlainRun $pos $discoveryText

# --- Remove the old "Built redistricting platform..." bullet entirely -----
$oldBuiltIdx = $P3 + 3
$full = $d.Paragraphs($oldBuiltIdx).Range
$r = $d.Range($full.Start, $full.End)
$r.Delete()

# "Achieved 87%...71%" has now shifted down to $oldBuiltIdx.
# --- Paragraph 6: 178% accuracy improvement --------------------------------
$pos = Clear-ParaContent $d.Paragraphs($oldBuiltIdx)
$t4 = $BULLET + " **"
$pos = Add-PlainRun $pos $t4
$pos = Add-BoldRun $pos "178%"
$t4b = "** accuracy improvement in racial classification algorithms"
$pos = Add-PlainRun $pos $t4b

Write-Output "Done"
